$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Commitment Date" column
$ws.Range("J1").Value = "Commitment Date"

# Rename the "Type" header to "Type * " (now a required field)
$ws.Range("C1").Value = "Type * "

# Commitment date for every existing row (2023-01-20, serial 44946)
$ws.Range("J2:J7").Value = 44946

# Apply the date number format to J2 first, then copy that formatting down
# so every cell in the column shares a single cell-format record.
$ws.Range("J2").NumberFormat = "mm-dd-yy"
$ws.Range("J2").Copy()
$ws.Range("J3:J7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Columns.Item(10).ColumnWidth = 14.86

$ws.Range("J3:J7").Select()
